$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.425663113594055
$ws.Range("B1").Value = 3.625338792800903
$ws.Range("C1").Value = 3.038676261901855
$ws.Range("D1").Value = 2.018335819244385
$ws.Range("E1").Value = 1.166325211524963
